# Input_POEntry.xlsx - "Added PO Approval / PO Receving" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 ("FileGroup") value text correction: SSR032022 -> SSR0302203
$ws.Range("B3").Value = "SSR0302203"

# Row 7 ("Purchase Qty") numeric value change: 3 -> 100
$ws.Range("B7").Value = 100

# Move/save the sheet's active selection to C1
$ws.Range("C1").Select()
